# Auto-generated edit script: applies verified cell-level diffs
# from the commit 'chore: update Sheets via scheduled runner'
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 19234.182
$ws.Range("I33").Value = 223
$ws.Range("J33").Value = 52503.75
$ws.Range("K33").Value = 223
$ws.Range("L33").Value = 52503.75
$ws.Range("M33").Value = 6
$ws.Range("N33").Value = -52961.75

# Row 101
$ws.Range("H101").Value = 1166.5
$ws.Range("I101").Value = 1168
$ws.Range("J101").Value = 1165
$ws.Range("K101").Value = 3504
$ws.Range("L101").Value = 3495
$ws.Range("M101").Value = -1882
$ws.Range("N101").Value = -6739

# Row 112
$ws.Range("H112").Value = 3798.2222
$ws.Range("J112").Value = 4160.5
$ws.Range("L112").Value = 12481.5
$ws.Range("N112").Value = -14697.5

# Row 135
$ws.Range("H135").Value = 966.5

# Row 138
$ws.Range("H138").Value = 5043.1113
$ws.Range("I138").Value = 1400
$ws.Range("J138").Value = 5498.5
$ws.Range("K138").Value = 4200
$ws.Range("L138").Value = 16495.5
$ws.Range("M138").Value = 940
$ws.Range("N138").Value = -26775.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("N32").ClearContents()
$ws.Range("H32").Value = 6761.3335
$ws.Range("I32").Value = 6761.3335
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 6761.3335
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -6474.3335

# Row 61
$ws.Range("H61").Value = 2215.75
$ws.Range("I61").Value = 913.1111
$ws.Range("J61").Value = 3281.5454
$ws.Range("K61").Value = 913.1111
$ws.Range("L61").Value = 3281.5454
$ws.Range("M61").Value = -701.1111
$ws.Range("N61").Value = -3705.5454

# Row 74
$ws.Range("N74").ClearContents()
$ws.Range("H74").Value = 1099.4
$ws.Range("I74").Value = 1099.4
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1099.4
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -225.4000000000001

# Row 77
$ws.Range("H77").Value = 1099.4
$ws.Range("I77").Value = 1099.4
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 5497
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -1129

# Row 136
$ws.Range("H136").Value = 2215.75
$ws.Range("I136").Value = 913.1111
$ws.Range("J136").Value = 3281.5454
$ws.Range("K136").Value = 2739.3333
$ws.Range("L136").Value = 9844.636200000001
$ws.Range("M136").Value = -189.3332999999998
$ws.Range("N136").Value = -14944.6362

# Row 138
$ws.Range("H138").Value = 84997.5
$ws.Range("J138").Value = 84997.5
$ws.Range("L138").Value = 84997.5
$ws.Range("N138").Value = -95277.5

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3824.75
$ws.Range("I86").Value = 2649.5
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 2649.5
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -1526.5
$ws.Range("N86").Value = -7246

# Row 89
$ws.Range("H89").Value = 3824.75
$ws.Range("I89").Value = 2649.5
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 13247.5
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -7631.5
$ws.Range("N89").Value = -36232

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 93.25
$ws.Range("I7").Value = 104.25
$ws.Range("J7").Value = 82.25
$ws.Range("K7").Value = 104.25
$ws.Range("L7").Value = 82.25
$ws.Range("M7").Value = 8.75
$ws.Range("N7").Value = -308.25

# Row 16
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0

# Row 31
$ws.Range("H31").Value = 3636.3635
$ws.Range("I31").Value = 1250
$ws.Range("K31").Value = 1250
$ws.Range("M31").Value = -955

# Row 34
$ws.Range("H34").Value = 3636.3635
$ws.Range("I34").Value = 1250
$ws.Range("K34").Value = 1250
$ws.Range("M34").Value = -1048

# Row 58
$ws.Range("H58").Value = 1000
$ws.Range("I58").Value = 1000
$ws.Range("K58").Value = 1000
$ws.Range("M58").Value = -797

# Row 113
$ws.Range("M113").ClearContents()
$ws.Range("N113").ClearContents()
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0

# Row 132
$ws.Range("H132").Value = 2083.7222
$ws.Range("I132").Value = 1558.5714
$ws.Range("K132").Value = 4675.7142
$ws.Range("M132").Value = -2145.7142

# Row 136
$ws.Range("H136").Value = 1000
$ws.Range("I136").Value = 1000
$ws.Range("K136").Value = 3000
$ws.Range("M136").Value = -450

# Row 141
$ws.Range("H141").Value = 99998.664
$ws.Range("J141").Value = 99998.664
$ws.Range("L141").Value = 99998.664
$ws.Range("N141").Value = -110358.664

$ws = $wb.Worksheets.Item("CUL")
# Row 21
$ws.Range("H21").Value = 874.5
$ws.Range("I21").Value = 1249
$ws.Range("K21").Value = 3747
$ws.Range("M21").Value = -3574

# Row 131
$ws.Range("H131").Value = 2594.8
$ws.Range("J131").Value = 2613.111
$ws.Range("L131").Value = 7839.333
$ws.Range("N131").Value = -17919.333

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1360.375
$ws.Range("I97").Value = 1357.8
$ws.Range("J97").Value = 1364.6666
$ws.Range("K97").Value = 1357.8
$ws.Range("L97").Value = 1364.6666
$ws.Range("M97").Value = -861.8
$ws.Range("N97").Value = -2356.6666

# Row 113
$ws.Range("N113").ClearContents()
$ws.Range("H113").Value = 3266.3333
$ws.Range("I113").Value = 3266.3333
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3266.3333
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1096.3333

# Row 122
$ws.Range("H122").Value = 11771.182
$ws.Range("I122").Value = 8948.299999999999
$ws.Range("K122").Value = 26844.9
$ws.Range("M122").Value = -24394.9

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 7950
$ws.Range("J122").Value = 7950
$ws.Range("L122").Value = 23850
$ws.Range("N122").Value = -28750

# Row 132
$ws.Range("H132").Value = 4154
$ws.Range("I132").Value = 3692.875
$ws.Range("K132").Value = 11078.625
$ws.Range("M132").Value = -8548.625

# Row 136
$ws.Range("H136").Value = 2588.5
$ws.Range("I136").Value = 2588.5
$ws.Range("K136").Value = 7765.5
$ws.Range("M136").Value = -5215.5

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("N122").ClearContents()
$ws.Range("H122").Value = 9999
$ws.Range("I122").Value = 9999
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 29997
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -27547

# Row 126
$ws.Range("M126").ClearContents()
$ws.Range("H126").Value = 7000.5
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 7000.5
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 21001.5
$ws.Range("N126").Value = -25941.5

# Row 132
$ws.Range("H132").Value = 2644.64
$ws.Range("I132").Value = 2374.652
$ws.Range("K132").Value = 7123.956
$ws.Range("M132").Value = -4593.956
